# Scheduled-runner price/profit refresh: updates currentAveragePrice* /
# LevePrice* / LeveProfit* columns (H:N) for a batch of leve rows across
# every crafting-job sheet, reflecting freshly pulled market data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 2238
$ws.Range("I88").Value = 1500
$ws.Range("J88").Value = 2976
$ws.Range("K88").Value = 1500
$ws.Range("L88").Value = 2976
$ws.Range("M88").Value = -1094
$ws.Range("N88").Value = -3788
$ws.Range("H91").Value = 2238
$ws.Range("I91").Value = 1500
$ws.Range("J91").Value = 2976
$ws.Range("K91").Value = 1500
$ws.Range("L91").Value = 2976
$ws.Range("M91").Value = -96
$ws.Range("N91").Value = -5784
$ws.Range("H113").Value = 3931.7334
$ws.Range("I113").Value = 2927.5
$ws.Range("J113").Value = 4296.909
$ws.Range("K113").Value = 2927.5
$ws.Range("L113").Value = 4296.909
$ws.Range("M113").Value = 326.5
$ws.Range("N113").Value = -10804.909
$ws.Range("H127").Value = 1127.2727
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 1127.2727
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 3381.8181
$ws.Range("M127").ClearContents()
$ws.Range("N127").Value = -13301.8181
$ws.Range("H137").Value = 1845.1951
$ws.Range("I137").Value = 1819.75
$ws.Range("J137").Value = 1900
$ws.Range("K137").Value = 5459.25
$ws.Range("L137").Value = 5700
$ws.Range("M137").Value = -2909.25
$ws.Range("N137").Value = -10800

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2671.6904
$ws.Range("I74").Value = 3316.577
$ws.Range("J74").Value = 1623.75
$ws.Range("K74").Value = 3316.577
$ws.Range("L74").Value = 1623.75
$ws.Range("M74").Value = -2442.577
$ws.Range("N74").Value = -3371.75
$ws.Range("H77").Value = 2671.6904
$ws.Range("I77").Value = 3316.577
$ws.Range("J77").Value = 1623.75
$ws.Range("K77").Value = 16582.885
$ws.Range("L77").Value = 8118.75
$ws.Range("M77").Value = -12214.885
$ws.Range("N77").Value = -16854.75
$ws.Range("H132").Value = 3943.7334
$ws.Range("I132").Value = 1245.5
$ws.Range("J132").Value = 10585.538
$ws.Range("K132").Value = 3736.5
$ws.Range("L132").Value = 31756.614
$ws.Range("M132").Value = -1206.5
$ws.Range("N132").Value = -36816.614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1364.75
$ws.Range("I86").Value = 1250
$ws.Range("J86").Value = 1479.5
$ws.Range("K86").Value = 1250
$ws.Range("L86").Value = 1479.5
$ws.Range("M86").Value = -127
$ws.Range("N86").Value = -3725.5
$ws.Range("H89").Value = 1364.75
$ws.Range("I89").Value = 1250
$ws.Range("J89").Value = 1479.5
$ws.Range("K89").Value = 6250
$ws.Range("L89").Value = 7397.5
$ws.Range("M89").Value = -634
$ws.Range("N89").Value = -18629.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6804822
$ws.Range("I31").Value = 1663.6487
$ws.Range("K31").Value = 1663.6487
$ws.Range("M31").Value = -1368.6487
$ws.Range("H34").Value = 6804822
$ws.Range("I34").Value = 1663.6487
$ws.Range("K34").Value = 1663.6487
$ws.Range("M34").Value = -1461.6487
$ws.Range("H80").Value = 50000
$ws.Range("J80").Value = 50000
$ws.Range("L80").Value = 50000
$ws.Range("N80").Value = -52246
$ws.Range("H83").Value = 50000
$ws.Range("J83").Value = 50000
$ws.Range("L83").Value = 150000
$ws.Range("N83").Value = -161232
$ws.Range("H132").Value = 2518.0356
$ws.Range("I132").Value = 1309
$ws.Range("J132").Value = 6145.143
$ws.Range("K132").Value = 3927
$ws.Range("L132").Value = 18435.429
$ws.Range("M132").Value = -1397
$ws.Range("N132").Value = -23495.429
$ws.Range("H134").Value = 2290.5
$ws.Range("I134").Value = 1288.7693
$ws.Range("J134").Value = 3158.6667
$ws.Range("K134").Value = 3866.3079
$ws.Range("L134").Value = 9476.000100000001
$ws.Range("M134").Value = -1331.3079
$ws.Range("N134").Value = -14546.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 509.65216
$ws.Range("I113").Value = 490.9091
$ws.Range("K113").Value = 1472.7273
$ws.Range("M113").Value = 697.2727
$ws.Range("H124").Value = 20854.545
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 20854.545
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 62563.63499999999
$ws.Range("M124").ClearContents()
$ws.Range("N124").Value = -72383.63499999999
$ws.Range("H131").Value = 371304.2
$ws.Range("I131").Value = 1429021.1
$ws.Range("J131").Value = 1103.25
$ws.Range("K131").Value = 4287063.300000001
$ws.Range("L131").Value = 3309.75
$ws.Range("M131").Value = -4282023.300000001
$ws.Range("N131").Value = -13389.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 10498.5
$ws.Range("J15").Value = 10498.5
$ws.Range("L15").Value = 10498.5
$ws.Range("N15").Value = -11074.5
$ws.Range("H63").Value = 11937.5
$ws.Range("J63").Value = 11937.5
$ws.Range("L63").Value = 11937.5
$ws.Range("N63").Value = -13309.5
$ws.Range("H66").Value = 11937.5
$ws.Range("J66").Value = 11937.5
$ws.Range("L66").Value = 35812.5
$ws.Range("N66").Value = -42676.5
$ws.Range("H81").Value = 10498.5
$ws.Range("J81").Value = 10498.5
$ws.Range("L81").Value = 10498.5
$ws.Range("N81").Value = -12494.5
$ws.Range("H84").Value = 10498.5
$ws.Range("J84").Value = 10498.5
$ws.Range("L84").Value = 31495.5
$ws.Range("N84").Value = -41479.5
$ws.Range("H102").Value = 2042992.4
$ws.Range("I102").Value = 2977708
$ws.Range("J102").Value = 3612.3635
$ws.Range("K102").Value = 2977708
$ws.Range("L102").Value = 3612.3635
$ws.Range("M102").Value = -2976086
$ws.Range("N102").Value = -6856.363499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4806.8125
$ws.Range("I61").Value = 3159.0833
$ws.Range("J61").Value = 9750
$ws.Range("K61").Value = 3159.0833
$ws.Range("L61").Value = 9750
$ws.Range("M61").Value = -2957.0833
$ws.Range("N61").Value = -10154
$ws.Range("H74").Value = 20000
$ws.Range("J74").Value = 20000
$ws.Range("L74").Value = 20000
$ws.Range("N74").Value = -21996
$ws.Range("H77").Value = 20000
$ws.Range("J77").Value = 20000
$ws.Range("L77").Value = 60000
$ws.Range("N77").Value = -69984
$ws.Range("H94").Value = 19865
$ws.Range("J94").Value = 19865
$ws.Range("L94").Value = 19865
$ws.Range("N94").Value = -21217
$ws.Range("H113").Value = 4806.8125
$ws.Range("I113").Value = 3159.0833
$ws.Range("J113").Value = 9750
$ws.Range("K113").Value = 3159.0833
$ws.Range("L113").Value = 9750
$ws.Range("M113").Value = -989.0832999999998
$ws.Range("N113").Value = -14090
$ws.Range("H122").Value = 4087.24
$ws.Range("I122").Value = 3865.611
$ws.Range("J122").Value = 4657.143
$ws.Range("K122").Value = 11596.833
$ws.Range("L122").Value = 13971.429
$ws.Range("M122").Value = -9146.832999999999
$ws.Range("N122").Value = -18871.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
